$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the version / last-tested-with text values on the sheet.
$ws.Range("A7").Value = "(Version: 1.0.1)"
$ws.Range("A8").Value = "(Last tested with: ReportServer 4.0.0-6053) "

# Move the active cell selection from A5 to A8 (matches the saved sheetView state).
$ws.Range("A8").Select()
